$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-06-04 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-05 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("73×92=6716", $true, $false, $false, $false, $false, $true, 1, $false, "43×19=817", 2) | Out-Null
$d.Content.Find.Execute("91×26=2366", $true, $false, $false, $false, $false, $true, 1, $false, "63×85=5355", 2) | Out-Null
$d.Content.Find.Execute("18×73=1314", $true, $false, $false, $false, $false, $true, 1, $false, "95×26=2470", 2) | Out-Null
$d.Content.Find.Execute("74×11=814", $true, $false, $false, $false, $false, $true, 1, $false, "56×77=4312", 2) | Out-Null
$d.Content.Find.Execute("26×75=1950", $true, $false, $false, $false, $false, $true, 1, $false, "68×50=3400", 2) | Out-Null
$d.Content.Find.Execute("38×34=1292", $true, $false, $false, $false, $false, $true, 1, $false, "37×15=555", 2) | Out-Null
$d.Content.Find.Execute("64×11=704", $true, $false, $false, $false, $false, $true, 1, $false, "95×13=1235", 2) | Out-Null
$d.Content.Find.Execute("26×31=806", $true, $false, $false, $false, $false, $true, 1, $false, "74×33=2442", 2) | Out-Null
$d.Content.Find.Execute("81×58=4698", $true, $false, $false, $false, $false, $true, 1, $false, "65×93=6045", 2) | Out-Null
$d.Content.Find.Execute("96×20=1920", $true, $false, $false, $false, $false, $true, 1, $false, "67×16=1072", 2) | Out-Null
$d.Content.Find.Execute("93×32=2976", $true, $false, $false, $false, $false, $true, 1, $false, "79×76=6004", 2) | Out-Null
$d.Content.Find.Execute("15×56=840", $true, $false, $false, $false, $false, $true, 1, $false, "26×71=1846", 2) | Out-Null
$d.Content.Find.Execute("71×60=4260", $true, $false, $false, $false, $false, $true, 1, $false, "56×57=3192", 2) | Out-Null
$d.Content.Find.Execute("46×91=4186", $true, $false, $false, $false, $false, $true, 1, $false, "15×84=1260", 2) | Out-Null
$d.Content.Find.Execute("92×93=8556", $true, $false, $false, $false, $false, $true, 1, $false, "79×17=1343", 2) | Out-Null
$d.Content.Find.Execute("66×69=4554", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=1344", 2) | Out-Null
$d.Content.Find.Execute("18×27=486", $true, $false, $false, $false, $false, $true, 1, $false, "54×34=1836", 2) | Out-Null
$d.Content.Find.Execute("59×17=1003", $true, $false, $false, $false, $false, $true, 1, $false, "82×43=3526", 2) | Out-Null
$d.Content.Find.Execute("84×39=3276", $true, $false, $false, $false, $false, $true, 1, $false, "11×64=704", 2) | Out-Null
$d.Content.Find.Execute("63×83=5229", $true, $false, $false, $false, $false, $true, 1, $false, "83×97=8051", 2) | Out-Null
$d.Content.Find.Execute("77×73=5621", $true, $false, $false, $false, $false, $true, 1, $false, "29×35=1015", 2) | Out-Null
$d.Content.Find.Execute("85×34=2890", $true, $false, $false, $false, $false, $true, 1, $false, "74×34=2516", 2) | Out-Null
$d.Content.Find.Execute("47×55=2585", $true, $false, $false, $false, $false, $true, 1, $false, "93×23=2139", 2) | Out-Null
$d.Content.Find.Execute("90×19=1710", $true, $false, $false, $false, $false, $true, 1, $false, "11×97=1067", 2) | Out-Null
$d.Content.Find.Execute("37×49=1813", $true, $false, $false, $false, $false, $true, 1, $false, "82×68=5576", 2) | Out-Null
